$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string value "diesel" via new row's Name cell
$ws.Range("A6").Value = "diesel"

# Numeric inputs for the new "diesel" row
$ws.Range("B6").Value = 600000
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 0.08
$ws.Range("H6").Value = 96.2

# Formulas consistent with the shared formulas already used in column G and I
$ws.Range("G6").Formula = "= (I6 + C6) /100*B6"
$ws.Range("I6").Formula = "=F6 / (1 - 1/(1 + F6) ^E6) *100"

# Update the active selection to match the new cell focus
$ws.Range("G6").Select()

$wb.Save()
